$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- OUTPUT table (rows 45-51): recalculated overflow split between the two
# must-have projects now that the overflow handling reads/accounts for the
# overflow amount instead of discarding it.
$ws.Range("B46").Value = 2.2000000000000002
$ws.Range("C46").Value = 5.8

$ws.Range("B47").Value = 4
$ws.Range("C47").Value = 4

$ws.Range("B48").Value = 2.8
$ws.Range("C48").Value = 2.8

$ws.Range("B49").Value = 2
$ws.Range("C49").Value = 2

$ws.Range("B50").Value = 4
$ws.Range("C50").Value = 4

$ws.Range("B51").Value = 4
$ws.Range("C51").Value = 4

# --- Overflow total (row 71) now reflects the overflow actually read in.
$ws.Range("B71").Value = 9.4

# --- Selection / scroll position left behind by the editing session.
[void]$ws.Range("B72").Select()

# --- Locale: default cell style is named "Standard" in the German build
# that produced this workbook (was "Normal").
$wb.Styles.Item("Normal").Name = "Standard"
